$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (activity values) and column C (frame, minutes) for rows 2-12
$ws.Range("B2").Value = 1595.1
$ws.Range("C2").Value = 13

$ws.Range("B3").Value = 158
$ws.Range("C3").Value = 13

$ws.Range("C4").Value = 13

$ws.Range("B5").Value = 1028
$ws.Range("C5").Value = 13

$ws.Range("B6").Value = 758
$ws.Range("C6").Value = 13

$ws.Range("B7").Value = 232
$ws.Range("C7").Value = 13

$ws.Range("C8").Value = 13

$ws.Range("B9").Value = 821
$ws.Range("C9").Value = 13

$ws.Range("B10").Value = 107
$ws.Range("C10").Value = 13

$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 13

$ws.Range("B12").Value = 74
$ws.Range("C12").Value = 13
